# Update the pathway comparison counting sheet with new experimental counts.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Difficidin
$ws.Range("B2").Value = 44
$ws.Range("C2").Value = 27
$ws.Range("E2").Value = 13

# Row 3 - Oocydin
$ws.Range("B3").Value = 37
$ws.Range("D3").Value = 16
$ws.Range("E3").Value = 10

# Remove the "S. lohii" note that was in H3
$ws.Range("H3").ClearContents()

# Row 4 - Bafilomycin
$ws.Range("B4").Value = 31

# Row 5 - Leupyrrin
$ws.Range("B5").Value = 19
$ws.Range("C5").Value = 12
$ws.Range("D5").Value = 6

# Row 6 - Tolaasin
$ws.Range("B6").Value = 20
$ws.Range("C6").Value = 17

# Row 7 - Anabaenopeptin
$ws.Range("B7").Value = 10
$ws.Range("C7").Value = 6
$ws.Range("D7").Value = 7

# Row 8 - Geldanamycin
$ws.Range("B8").Value = 24
$ws.Range("C8").Value = 24

# Row 9 - Oxazolomycin
$ws.Range("B9").Value = 43
$ws.Range("C9").Value = 28

# Update the selected cell as it was left in the authored workbook
$ws.Range("C17").Select()
